# Mexico Liga MX Femenil - base update (21-04-2024)
#
# 1) Seven pairs of adjacent match rows had their fixture details (everything
#    except the running index in column A and the constant Div/Div Original
#    Name/Date columns C:E, which were already identical within each pair)
#    swapped between the two rows.
# 2) Two brand-new fixtures were inserted right before the old last row
#    (296), pushing the previous rows 296/297 down to 298/299 and bumping
#    their running index (column A) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    # Column B (id) swap
    $cb1 = $ws.Range("B$rowA")
    $cb2 = $ws.Range("B$rowB")
    $vb1 = $cb1.Value2
    $vb2 = $cb2.Value2
    $cb1.Value2 = $vb2
    $cb2.Value2 = $vb1

    # Columns F:AC (everything after the constant Div/Div Original Name/Date
    # block) swap
    $cf1 = $ws.Range("F" + $rowA + ":AC" + $rowA)
    $cf2 = $ws.Range("F" + $rowB + ":AC" + $rowB)
    $vf1 = $cf1.Value2
    $vf2 = $cf2.Value2
    $cf1.Value2 = $vf2
    $cf2.Value2 = $vf1
}

Swap-Rows 28 29
Swap-Rows 71 72
Swap-Rows 133 134
Swap-Rows 149 150
Swap-Rows 232 233
Swap-Rows 245 246
Swap-Rows 271 272

# Insert two new rows before the old row 296 (xlShiftDown, format from
# above so the running-index/date style carries over).
$ws.Rows("296:297").Insert(-4121, 0)

# The insert copies the format of row 295 which (per engine quirk) does not
# keep the bordered/bold style used by column A, nor the custom date format
# on column E -- repair both by pasting the format down from the row that
# used to be 296 (now pushed to 298, which still carries the original
# style).
$ws.Range("A298").Copy() | Out-Null
$ws.Range("A296:A297").PasteSpecial(-4122) | Out-Null
$ws.Range("E298").Copy() | Out-Null
$ws.Range("E296:E297").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New row 296 data
$data296 = New-Object 'object[,]' 1,29
$vals296 = @(294, 7645828, "Mexico Liga MX Femenil", "Mexico Liga MX Femenil", 45401.83333333334, "Queretaro Women", "Leon Women", 3, 2, "H", 2.35, 3.6, 2.45, 2.8, 3.6, 2.05, 0.25, 1.9, 1.9, 2.5, 1.75, 1.95, 1.8, -1, -1, 0.8999999999999999, -1, 0.75, -1)
for ($i = 0; $i -lt 29; $i++) { $data296[0,$i] = $vals296[$i] }
$ws.Range("A296:AC296").Value2 = $data296

# New row 297 data
$data297 = New-Object 'object[,]' 1,29
$vals297 = @(295, 7645735, "Mexico Liga MX Femenil", "Mexico Liga MX Femenil", 45401.91666666666, "Atlas Women", "Tigres UANL Women", 0, 1, "A", 12, 8, 1.142, 12, 8, 1.142, 2.5, 1.9, 1.9, 4.25, 1.9, 1.9, -1, -1, 0.1419999999999999, 0.8999999999999999, -1, -1, 0.8999999999999999)
for ($i = 0; $i -lt 29; $i++) { $data297[0,$i] = $vals297[$i] }
$ws.Range("A297:AC297").Value2 = $data297

# The two rows that used to be 296/297 are now 298/299 (Excel already moved
# their contents down); only their running index (column A) needs bumping.
$ws.Range("A298").Value2 = 296
$ws.Range("A299").Value2 = 297
